$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $val) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "26.326.18"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "1.611.01"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  -0.16%  "
Set-TextValue $ws.Range("D5") "213.24"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("E9").Value = "  +0.25%  "
Set-TextValue $ws.Range("D10") "18.53"
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "1.834.26"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "1.601.80"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("E14").Value = "  +0.56%  "
Set-TextValue $ws.Range("D15") "0.517"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").Value = "26.297.93"
$ws.Range("E16").Value = "  +0.91%  "
Set-TextValue $ws.Range("D17") "62.39"
$ws.Range("E17").Value = "  +3.38%  "
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("E19").Value = "  -0.17%  "
Set-TextValue $ws.Range("D20") "202.62"
$ws.Range("E20").Value = "  +0.46%  "
Set-TextValue $ws.Range("D21") "4.28"
$ws.Range("E21").Value = "  +1.31%  "
Set-TextValue $ws.Range("D22") "9.36"
$ws.Range("E22").Value = "  +0.95%  "
$ws.Range("E23").Value = "  +1.03%  "
$ws.Range("E24").Value = "  +1.77%  "
Set-TextValue $ws.Range("D25") "143.54"
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("E27").Value = "  -0.60%  "
Set-TextValue $ws.Range("D28") "15.26"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("E29").Value = "  +2.52%  "
$ws.Range("E30").Value = "  +5.40%  "
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("E32").Value = "  +3.07%  "
$ws.Range("E33").Value = "  -0.21%  "
Set-TextValue $ws.Range("D34") "1.50"
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("E35").Value = "  +1.08%  "
$ws.Range("D36").Value = "1.163.66"
$ws.Range("E36").Value = "  +3.43%  "
$ws.Range("E37").Value = "  +1.81%  "
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("E41").Value = "  +1.16%  "
Set-TextValue $ws.Range("D42") "5.39"
$ws.Range("E42").Value = "  +4.46%  "
Set-TextValue $ws.Range("D43") "0.784"
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").Value = "1.745.14"
$ws.Range("E44").Value = "  +0.36%  "
Set-TextValue $ws.Range("D45") "92.53"
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D46") "1.54"
$ws.Range("E46").Value = "  +1.65%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0105"
$ws.Range("E47").Value = "  +13.39%  "
Set-TextValue $ws.Range("D48") "53.97"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("E49").Value = "  +0.42%  "
Set-TextValue $ws.Range("D50") "0.408"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("E51").Value = "  -0.26%  "
